$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2857142857142857
$ws.Range("C2").Value = 0.5555555555555556
$ws.Range("D2").Value = 0.5294117647058824
$ws.Range("E2").Value = 0.6666666666666667
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.5882352941176471
$ws.Range("H2").Value = 0.7894736842105263
$ws.Range("I2").Value = 0.85
$ws.Range("J2").Value = 0.6470588235294117
$ws.Range("K2").Value = 0.7619047619047619
$ws.Range("L2").Value = 0.6111111111111112
$ws.Range("M2").Value = 0.85
$ws.Range("N2").Value = 0.5
$ws.Range("O2").Value = 0.7222222222222222
$ws.Range("P2").Value = 0.7
$ws.Range("A3").Value = 0.2857142857142857
$ws.Range("C3").Value = 0.5555555555555556
$ws.Range("D3").Value = 0.4375
$ws.Range("E3").Value = 0.6666666666666667
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.6666666666666667
$ws.Range("H3").Value = 0.7894736842105263
$ws.Range("I3").Value = 0.85
$ws.Range("J3").Value = 0.6470588235294117
$ws.Range("K3").Value = 0.7619047619047619
$ws.Range("L3").Value = 0.5294117647058824
$ws.Range("M3").Value = 0.85
$ws.Range("N3").Value = 0.5
$ws.Range("O3").Value = 0.7222222222222222
$ws.Range("P3").Value = 0.7
$ws.Range("A4").Value = 0.5555555555555556
$ws.Range("B4").Value = 0.5555555555555556
$ws.Range("D4").Value = 0.5789473684210527
$ws.Range("E4").Value = 0.4705882352941176
$ws.Range("F4").Value = 0.375
$ws.Range("G4").Value = 0.4705882352941176
$ws.Range("H4").Value = 0.4375
$ws.Range("I4").Value = 0.5294117647058824
$ws.Range("J4").Value = 0.4375
$ws.Range("K4").Value = 0.5263157894736843
$ws.Range("L4").Value = 0.4117647058823529
$ws.Range("M4").Value = 0.75
$ws.Range("N4").Value = 0.4736842105263158
$ws.Range("O4").Value = 0.5294117647058824
$ws.Range("P4").Value = 0.4444444444444444
$ws.Range("A5").Value = 0.5294117647058824
$ws.Range("B5").Value = 0.4375
$ws.Range("C5").Value = 0.5789473684210527
$ws.Range("E5").Value = 0.5294117647058824
$ws.Range("F5").Value = 0.3333333333333334
$ws.Range("G5").Value = 0.5294117647058824
$ws.Range("H5").Value = 0.8
$ws.Range("I5").Value = 0.736842105263158
$ws.Range("J5").Value = 0.5
$ws.Range("K5").Value = 0.65
$ws.Range("L5").Value = 0.2666666666666667
$ws.Range("M5").Value = 0.8
$ws.Range("N5").Value = 0.25
$ws.Range("O5").Value = 0.5882352941176471
$ws.Range("P5").Value = 0.5789473684210527
$ws.Range("A6").Value = 0.6666666666666667
$ws.Range("B6").Value = 0.6666666666666667
$ws.Range("C6").Value = 0.4705882352941176
$ws.Range("D6").Value = 0.5294117647058824
$ws.Range("F6").Value = 0.2857142857142857
$ws.Range("G6").Value = 0.2857142857142857
$ws.Range("H6").Value = 0.6470588235294117
$ws.Range("I6").Value = 0.5625
$ws.Range("J6").Value = 0.2307692307692307
$ws.Range("K6").Value = 0.2666666666666667
$ws.Range("L6").Value = 0.3333333333333334
$ws.Range("M6").Value = 0.4666666666666667
$ws.Range("N6").Value = 0.4117647058823529
$ws.Range("O6").Value = 0.3571428571428571
$ws.Range("P6").Value = 0.2666666666666667
$ws.Range("A7").Value = 0.5
$ws.Range("B7").Value = 0.5
$ws.Range("C7").Value = 0.375
$ws.Range("D7").Value = 0.3333333333333334
$ws.Range("E7").Value = 0.2857142857142857
$ws.Range("G7").Value = 0.2857142857142857
$ws.Range("H7").Value = 0.6470588235294117
$ws.Range("I7").Value = 0.5625
$ws.Range("J7").Value = 0.2307692307692307
$ws.Range("K7").Value = 0.4705882352941176
$ws.Range("L7").Value = 0.2142857142857143
$ws.Range("M7").Value = 0.6470588235294117
$ws.Range("N7").Value = 0.3125
$ws.Range("O7").Value = 0.3571428571428571
$ws.Range("P7").Value = 0.375
$ws.Range("A8").Value = 0.5882352941176471
$ws.Range("B8").Value = 0.6666666666666667
$ws.Range("C8").Value = 0.4705882352941176
$ws.Range("D8").Value = 0.5294117647058824
$ws.Range("E8").Value = 0.2857142857142857
$ws.Range("F8").Value = 0.2857142857142857
$ws.Range("H8").Value = 0.5625
$ws.Range("I8").Value = 0.4666666666666667
$ws.Range("J8").Value = 0.08333333333333337
$ws.Range("K8").Value = 0.4705882352941176
$ws.Range("L8").Value = 0.3333333333333334
$ws.Range("M8").Value = 0.5625
$ws.Range("N8").Value = 0.4117647058823529
$ws.Range("O8").Value = 0.2307692307692307
$ws.Range("P8").Value = 0.4705882352941176
$ws.Range("A9").Value = 0.7894736842105263
$ws.Range("B9").Value = 0.7894736842105263
$ws.Range("C9").Value = 0.4375
$ws.Range("D9").Value = 0.8
$ws.Range("E9").Value = 0.6470588235294117
$ws.Range("F9").Value = 0.6470588235294117
$ws.Range("G9").Value = 0.5625
$ws.Range("I9").Value = 0.1666666666666666
$ws.Range("J9").Value = 0.5333333333333333
$ws.Range("K9").Value = 0.6111111111111112
$ws.Range("M9").Value = 0.7058823529411764
$ws.Range("N9").Value = 0.7
$ws.Range("O9").Value = 0.625
$ws.Range("A10").Value = 0.85
$ws.Range("B10").Value = 0.85
$ws.Range("C10").Value = 0.5294117647058824
$ws.Range("D10").Value = 0.736842105263158
$ws.Range("E10").Value = 0.5625
$ws.Range("F10").Value = 0.5625
$ws.Range("G10").Value = 0.4666666666666667
$ws.Range("H10").Value = 0.1666666666666666
$ws.Range("J10").Value = 0.4285714285714286
$ws.Range("K10").Value = 0.5294117647058824
$ws.Range("L10").Value = 0.5882352941176471
$ws.Range("M10").Value = 0.625
$ws.Range("N10").Value = 0.631578947368421
$ws.Range("O10").Value = 0.5333333333333333
$ws.Range("P10").Value = 0.6842105263157895
$ws.Range("A11").Value = 0.6470588235294117
$ws.Range("B11").Value = 0.6470588235294117
$ws.Range("C11").Value = 0.4375
$ws.Range("D11").Value = 0.5
$ws.Range("E11").Value = 0.2307692307692307
$ws.Range("F11").Value = 0.2307692307692307
$ws.Range("G11").Value = 0.08333333333333337
$ws.Range("H11").Value = 0.5333333333333333
$ws.Range("I11").Value = 0.4285714285714286
$ws.Range("K11").Value = 0.4375
$ws.Range("L11").Value = 0.2857142857142857
$ws.Range("M11").Value = 0.5333333333333333
$ws.Range("N11").Value = 0.375
$ws.Range("O11").Value = 0.1666666666666666
$ws.Range("P11").Value = 0.4375
$ws.Range("A12").Value = 0.7619047619047619
$ws.Range("B12").Value = 0.7619047619047619
$ws.Range("C12").Value = 0.5263157894736843
$ws.Range("D12").Value = 0.65
$ws.Range("E12").Value = 0.2666666666666667
$ws.Range("F12").Value = 0.4705882352941176
$ws.Range("G12").Value = 0.4705882352941176
$ws.Range("H12").Value = 0.6111111111111112
$ws.Range("I12").Value = 0.5294117647058824
$ws.Range("J12").Value = 0.4375
$ws.Range("L12").Value = 0.5
$ws.Range("M12").Value = 0.3333333333333334
$ws.Range("N12").Value = 0.4736842105263158
$ws.Range("O12").Value = 0.4375
$ws.Range("P12").Value = 0.25
$ws.Range("A13").Value = 0.6111111111111112
$ws.Range("B13").Value = 0.5294117647058824
$ws.Range("C13").Value = 0.4117647058823529
$ws.Range("D13").Value = 0.2666666666666667
$ws.Range("E13").Value = 0.3333333333333334
$ws.Range("F13").Value = 0.2142857142857143
$ws.Range("G13").Value = 0.3333333333333334
$ws.Range("I13").Value = 0.5882352941176471
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("K13").Value = 0.5
$ws.Range("N13").Value = 0.3529411764705882
$ws.Range("O13").Value = 0.4
$ws.Range("P13").Value = 0.4117647058823529
$ws.Range("A14").Value = 0.85
$ws.Range("B14").Value = 0.85
$ws.Range("C14").Value = 0.75
$ws.Range("D14").Value = 0.8
$ws.Range("E14").Value = 0.4666666666666667
$ws.Range("F14").Value = 0.6470588235294117
$ws.Range("G14").Value = 0.5625
$ws.Range("H14").Value = 0.7058823529411764
$ws.Range("I14").Value = 0.625
$ws.Range("J14").Value = 0.5333333333333333
$ws.Range("K14").Value = 0.3333333333333334
$ws.Range("N14").Value = 0.631578947368421
$ws.Range("O14").Value = 0.4285714285714286
$ws.Range("P14").Value = 0.4375
$ws.Range("A15").Value = 0.5
$ws.Range("B15").Value = 0.5
$ws.Range("C15").Value = 0.4736842105263158
$ws.Range("D15").Value = 0.25
$ws.Range("E15").Value = 0.4117647058823529
$ws.Range("F15").Value = 0.3125
$ws.Range("G15").Value = 0.4117647058823529
$ws.Range("H15").Value = 0.7
$ws.Range("I15").Value = 0.631578947368421
$ws.Range("J15").Value = 0.375
$ws.Range("K15").Value = 0.4736842105263158
$ws.Range("L15").Value = 0.3529411764705882
$ws.Range("M15").Value = 0.631578947368421
$ws.Range("O15").Value = 0.375
$ws.Range("P15").Value = 0.3888888888888888
$ws.Range("A16").Value = 0.7222222222222222
$ws.Range("B16").Value = 0.7222222222222222
$ws.Range("C16").Value = 0.5294117647058824
$ws.Range("D16").Value = 0.5882352941176471
$ws.Range("E16").Value = 0.3571428571428571
$ws.Range("F16").Value = 0.3571428571428571
$ws.Range("G16").Value = 0.2307692307692307
$ws.Range("H16").Value = 0.625
$ws.Range("I16").Value = 0.5333333333333333
$ws.Range("J16").Value = 0.1666666666666666
$ws.Range("K16").Value = 0.4375
$ws.Range("L16").Value = 0.4
$ws.Range("M16").Value = 0.4285714285714286
$ws.Range("N16").Value = 0.375
$ws.Range("A17").Value = 0.7
$ws.Range("B17").Value = 0.7
$ws.Range("C17").Value = 0.4444444444444444
$ws.Range("D17").Value = 0.5789473684210527
$ws.Range("E17").Value = 0.2666666666666667
$ws.Range("F17").Value = 0.375
$ws.Range("G17").Value = 0.4705882352941176
$ws.Range("I17").Value = 0.6842105263157895
$ws.Range("J17").Value = 0.4375
$ws.Range("K17").Value = 0.25
$ws.Range("L17").Value = 0.4117647058823529
$ws.Range("M17").Value = 0.4375
$ws.Range("N17").Value = 0.3888888888888888
